# Insert a new weekly price record above the current row 350 ("Cebollín"
# feria Lagunitas de Puerto Montt data set). This shifts the existing
# rows 350-418 down to 351-419 (the former row 418 becomes row 419) and
# the new row 350 is filled with the new week's price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 350, pushing rows 350:418 down to 351:419.
$ws.Rows.Item(350).Insert()

# Populate the newly inserted row 350 with the new record's values.
$ws.Cells.Item(350, 1).Value = 4
$ws.Cells.Item(350, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(350, 3).Value = "Los Lagos"
$ws.Cells.Item(350, 4).Value = 44995
$ws.Cells.Item(350, 5).Value = 10
$ws.Cells.Item(350, 6).Value = 100112037
$ws.Cells.Item(350, 7).Value = "Cebollín"
$ws.Cells.Item(350, 8).Value = "Sin especificar"
$ws.Cells.Item(350, 9).Value = "Primera"
$ws.Cells.Item(350, 10).Value = 160
$ws.Cells.Item(350, 11).Value = 6500
$ws.Cells.Item(350, 12).Value = 7000
$ws.Cells.Item(350, 13).Value = 6750
$ws.Cells.Item(350, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(350, 15).Value = "Región Metropolitana"
$ws.Cells.Item(350, 16).Value = 188
$ws.Cells.Item(350, 17).Value = 36
$ws.Cells.Item(350, 18).Value = "Hortaliza"
